# Auto-generated edit script applying numeric cell updates per the commit diff.
# For each changed cell: set new value, or clear the cell if the diff removes it entirely.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 80
$ws.Range("I8").Value = 80
$ws.Range("K8").Value = 240
$ws.Range("M8").Value = -101

$ws.Range("H41").Value = 468.2857
$ws.Range("I41").Value = 395.25
$ws.Range("K41").Value = 395.25
$ws.Range("M41").Value = 44.75

$ws.Range("H43").Value = 800
$ws.Range("J43").Value = 800
$ws.Range("L43").Value = 800
$ws.Range("N43").Value = -938

$ws.Range("H116").Value = 4025
$ws.Range("I116").Value = 4025
$ws.Range("K116").Value = 4025
$ws.Range("M116").Value = -583

$ws.Range("H138").Value = 3923.15
$ws.Range("I138").Value = 1760.6666
$ws.Range("K138").Value = 5281.9998
$ws.Range("M138").Value = -141.9997999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1733.4
$ws.Range("I2").Value = 1733.4
$ws.Range("K2").Value = 1733.4
$ws.Range("M2").Value = -1620.4

$ws.Range("H4").Value = 117.5
$ws.Range("I4").Value = 90
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 90
$ws.Range("L4").Value = 200
$ws.Range("M4").Value = 26
$ws.Range("N4").Value = -432

$ws.Range("H32").Value = 7769.3335
$ws.Range("I32").Value = 7354.294
$ws.Range("K32").Value = 7354.294
$ws.Range("M32").Value = -7067.294

$ws.Range("H45").Value = 2216.5
$ws.Range("I45").Value = 2250
$ws.Range("K45").Value = 2250
$ws.Range("M45").Value = -1873

$ws.Range("H61").Value = 2642.6
$ws.Range("I61").Value = 2642.6
$ws.Range("K61").Value = 2642.6
$ws.Range("M61").Value = -2430.6

$ws.Range("H74").Value = 18444.615
$ws.Range("I74").Value = 17485.9
$ws.Range("J74").Value = 21640.334
$ws.Range("K74").Value = 17485.9
$ws.Range("L74").Value = 21640.334
$ws.Range("M74").Value = -16611.9
$ws.Range("N74").Value = -23388.334

$ws.Range("H77").Value = 18444.615
$ws.Range("I77").Value = 17485.9
$ws.Range("J77").Value = 21640.334
$ws.Range("K77").Value = 87429.5
$ws.Range("L77").Value = 108201.67
$ws.Range("M77").Value = -83061.5
$ws.Range("N77").Value = -116937.67

$ws.Range("H116").Value = 1733.4
$ws.Range("I116").Value = 1733.4
$ws.Range("K116").Value = 1733.4
$ws.Range("M116").Value = 560.5999999999999

$ws.Range("H132").Value = 2462.3333
$ws.Range("I132").Value = 2508.3333
$ws.Range("J132").Value = 2416.3333
$ws.Range("K132").Value = 7524.999899999999
$ws.Range("L132").Value = 7248.999899999999
$ws.Range("M132").Value = -4994.999899999999
$ws.Range("N132").Value = -12308.9999

$ws.Range("H136").Value = 2642.6
$ws.Range("I136").Value = 2642.6
$ws.Range("K136").Value = 7927.799999999999
$ws.Range("M136").Value = -5377.799999999999

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").ClearContents()
$ws.Range("N139").Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1733.4
$ws.Range("I3").Value = 1733.4
$ws.Range("K3").Value = 1733.4
$ws.Range("M3").Value = -1619.4

$ws.Range("H8").Value = 500
$ws.Range("J8").Value = 500
$ws.Range("L8").Value = 500
$ws.Range("N8").Value = -780

$ws.Range("H80").Value = 1670.8572
$ws.Range("I80").Value = 1480
$ws.Range("J80").Value = 1702.6666
$ws.Range("K80").Value = 1480
$ws.Range("L80").Value = 1702.6666
$ws.Range("M80").Value = -482
$ws.Range("N80").Value = -3698.6666

$ws.Range("H83").Value = 1670.8572
$ws.Range("I83").Value = 1480
$ws.Range("J83").Value = 1702.6666
$ws.Range("K83").Value = 7400
$ws.Range("L83").Value = 8513.333000000001
$ws.Range("M83").Value = -2408
$ws.Range("N83").Value = -18497.333

$ws.Range("H99").Value = 26000.25
$ws.Range("I99").Value = 26000.25
$ws.Range("K99").Value = 26000.25
$ws.Range("M99").Value = -24502.25

$ws.Range("H107").Value = 2466.1667
$ws.Range("I107").Value = 2469.8
$ws.Range("K107").Value = 2469.8
$ws.Range("M107").Value = -549.8000000000002

$ws.Range("H134").Value = 2781.0908
$ws.Range("I134").Value = 2781.0908
$ws.Range("K134").Value = 8343.2724
$ws.Range("M134").Value = -5808.2724

$ws.Range("H135").Value = 65283.223
$ws.Range("J135").Value = 65283.223
$ws.Range("L135").Value = 65283.223
$ws.Range("N135").Value = -75423.223

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 199
$ws.Range("I7").Value = 179
$ws.Range("K7").Value = 179
$ws.Range("M7").Value = -66

$ws.Range("H16").Value = 3539.5
$ws.Range("I16").Value = 2219
$ws.Range("K16").Value = 2219
$ws.Range("M16").Value = -1932

$ws.Range("H22").Value = 316.66666
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws.Range("H31").Value = 2349.1667
$ws.Range("I31").Value = 1993.75
$ws.Range("J31").Value = 3060
$ws.Range("K31").Value = 1993.75
$ws.Range("L31").Value = 3060
$ws.Range("M31").Value = -1698.75
$ws.Range("N31").Value = -3650

$ws.Range("H34").Value = 2349.1667
$ws.Range("I34").Value = 1993.75
$ws.Range("J34").Value = 3060
$ws.Range("K34").Value = 1993.75
$ws.Range("L34").Value = 3060
$ws.Range("M34").Value = -1791.75
$ws.Range("N34").Value = -3464

$ws.Range("H58").Value = 4199.2856
$ws.Range("I58").Value = 2365
$ws.Range("J58").Value = 5575
$ws.Range("K58").Value = 2365
$ws.Range("L58").Value = 5575
$ws.Range("M58").Value = -2162
$ws.Range("N58").Value = -5981

$ws.Range("H99").Value = 3000
$ws.Range("I99").Value = 3000
$ws.Range("K99").Value = 3000
$ws.Range("M99").Value = -1502

$ws.Range("H105").Value = 3332.3333
$ws.Range("I105").Value = 3332.3333
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 3332.3333
$ws.Range("L105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = -1585.3333

$ws.Range("H113").Value = 3539.5
$ws.Range("I113").Value = 2219
$ws.Range("K113").Value = 2219
$ws.Range("M113").Value = -49

$ws.Range("H122").Value = 5499
$ws.Range("I122").Value = 4798.6
$ws.Range("J122").Value = 7250
$ws.Range("K122").Value = 14395.8
$ws.Range("L122").Value = 21750
$ws.Range("M122").Value = -11945.8
$ws.Range("N122").Value = -26650

$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -6530

$ws.Range("H132").Value = 4562
$ws.Range("I132").Value = 3933
$ws.Range("J132").Value = 4939.4
$ws.Range("K132").Value = 11799
$ws.Range("L132").Value = 14818.2
$ws.Range("M132").Value = -9269
$ws.Range("N132").Value = -19878.2

$ws.Range("H134").Value = 2194.8333
$ws.Range("I134").Value = 2227
$ws.Range("J134").Value = 2178.75
$ws.Range("K134").Value = 6681
$ws.Range("L134").Value = 6536.25
$ws.Range("M134").Value = -4146
$ws.Range("N134").Value = -11606.25

$ws.Range("H136").Value = 4199.2856
$ws.Range("I136").Value = 2365
$ws.Range("J136").Value = 5575
$ws.Range("K136").Value = 7095
$ws.Range("L136").Value = 16725
$ws.Range("M136").Value = -4545
$ws.Range("N136").Value = -21825

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2732.4167
$ws.Range("I3").Value = 2732.4167
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 8197.250100000001
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -8085.250100000001

$ws.Range("H121").Value = 2474.3333
$ws.Range("I121").Value = 1089.6666
$ws.Range("K121").Value = 3268.9998
$ws.Range("M121").Value = -1958.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 999.6667
$ws.Range("I102").Value = 999.6667
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 999.6667
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = 622.3333

$ws.Range("H107").Value = 2545.125
$ws.Range("I107").Value = 3152
$ws.Range("J107").Value = 1533.6666
$ws.Range("K107").Value = 3152
$ws.Range("L107").Value = 1533.6666
$ws.Range("M107").Value = -1232
$ws.Range("N107").Value = -5373.6666

$ws.Range("H122").Value = 3150.5
$ws.Range("I122").Value = 2700.6667
$ws.Range("K122").Value = 8102.000100000001
$ws.Range("M122").Value = -5652.000100000001

$ws.Range("H132").Value = 3869.7273
$ws.Range("I132").Value = 3345.5
$ws.Range("K132").Value = 10036.5
$ws.Range("M132").Value = -7506.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3772.25
$ws.Range("I122").Value = 3696.3333
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 11088.9999
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -8638.999899999999
$ws.Range("N122").Value = -16900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4064.7778
$ws.Range("I132").Value = 4446.5
$ws.Range("J132").Value = 3759.4
$ws.Range("K132").Value = 13339.5
$ws.Range("L132").Value = 11278.2
$ws.Range("M132").Value = -10809.5
$ws.Range("N132").Value = -16338.2

